$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.144.35"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "1.779.32"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "31.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.28%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.291"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0690"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0945"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "2.035.62"
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "1.791.64"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").Value = "34.112.82"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.621"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.88"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("E19").Value = "  +3.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("E21").Value = "  +3.52%  "
$ws.Range("E22").Value = "  +0.19%  "
$ws.Range("E23").Value = "  +1.43%  "
$ws.Range("E24").Value = "  -0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0520"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.13%  "
$ws.Range("E33").Value = "  +5.82%  "
$ws.Range("D35").Value = "1.442.12"
$ws.Range("E35").Value = "  +3.56%  "
$ws.Range("E36").Value = "  +3.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.12%  "
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.22"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.920"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  -0.87%  "
$ws.Range("E45").Value = "  +0.75%  "
$ws.Range("E46").Value = "  +3.58%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").Value = "1.938.03"
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.43%  "
$ws.Range("E51").Value = "  +0.22%  "
